{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same change as the target diff:\n//  1. Splits the \"Because the first catch...\" run so that\n//     \"ArithmeticException\" is wrapped in proofErr spellStart/spellEnd.\n//  2. Splits the \"B,D,A,C\" run so that \",D,A,C\" is wrapped in\n//     proofErr gramStart/gramEnd (keeping \"B\" as its own run).\n//  3. Splits \"In .java files attached\" into \"In .java\" (wrapped in\n//     proofErr gramStart/gramEnd) + \" files attached\", moves the\n//     \"_GoBack\" bookmark into its own (new, empty) paragraph right\n//     after it, and removes the trailing \"       1.\" / \"       2.\"\n//     paragraphs entirely.\n\nconst FLAT_OPC_NS = \"http://schemas.microsoft.com/office/2006/xmlPackage\";\nconst WORD_MAIN_CT =\n  \"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\";\n\n// Helper: wrap a fragment of <w:body> children in the flat-OPC \"package\"\n// envelope that Range.insertOoxml expects.\nfunction flatOpc(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    `<pkg:package xmlns:pkg=\"${FLAT_OPC_NS}\">` +\n    `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"${WORD_MAIN_CT}\">` +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body>${bodyInnerXml}</w:body>` +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// pPr shared by the two \"ListParagraph\" bullet items we touch below.\nconst LIST_PPR =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>';\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three paragraphs we need purely by their current text, so the\n// script is resilient to the exact paragraph index.\nlet idxException = -1;\nlet idxBdac = -1;\nlet idxJava = -1;\nlet idxOne = -1;\nlet idxTwo = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"second catch statement with ArithmeticException\") !== -1) {\n    idxException = i;\n  } else if (t === \"B,D,A,C\") {\n    idxBdac = i;\n  } else if (t === \"In .java files attached\") {\n    idxJava = i;\n  } else if (t.trim() === \"1.\") {\n    idxOne = i;\n  } else if (t.trim() === \"2.\") {\n    idxTwo = i;\n  }\n}\n\nif (idxException === -1 || idxBdac === -1 || idxJava === -1) {\n  throw new Error(\n    \"Could not locate expected paragraphs (exception=\" +\n      idxException +\n      \", bdac=\" +\n      idxBdac +\n      \", java=\" +\n      idxJava +\n      \")\"\n  );\n}\n\n// 1) \"Because the first catch statement ... ArithmeticException will never\n//    be reached.\" -> split into 3 runs with a spellStart/spellEnd proofErr\n//    pair wrapping \"ArithmeticException\".\nconst exceptionOoxml = flatOpc(\n  \"<w:p>\" +\n    LIST_PPR +\n    '<w:r><w:t xml:space=\"preserve\">Because the first catch statement catches the type Exception, all exceptions will be caught. Because of this, the second catch statement with </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>ArithmeticException</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> will never be reached.</w:t></w:r>' +\n    \"</w:p>\"\n);\nparagraphs.items[idxException]\n  .getRange()\n  .insertOoxml(exceptionOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \"B,D,A,C\" -> \"B\" + proofErr gramStart/gramEnd wrapping \",D,A,C\".\nconst bdacOoxml = flatOpc(\n  \"<w:p>\" +\n    LIST_PPR +\n    \"<w:r><w:t>B</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>,D,A,C</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"</w:p>\"\n);\nparagraphs.items[idxBdac]\n  .getRange()\n  .insertOoxml(bdacOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) \"In .java files attached\" (+ trailing _GoBack bookmark, + the\n//    \"       1.\" / \"       2.\" paragraphs that follow) all become:\n//      <p> proofErr gramStart, \"In .java\", proofErr gramEnd, \" files\n//          attached\" </p>\n//      <p> bookmarkStart/_GoBack, bookmarkEnd </p>\n//    i.e. the bookmark moves onto its own empty paragraph and the\n//    \"1.\"/\"2.\" paragraphs disappear.\nparagraphs.load(\"items/text\");\nawait context.sync();\n// Re-resolve indices (they have not shifted from the two edits above,\n// since those were in-place replacements of the same paragraph count).\nidxJava = -1;\nidxOne = -1;\nidxTwo = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"In .java files attached\") {\n    idxJava = i;\n  } else if (t.trim() === \"1.\") {\n    idxOne = i;\n  } else if (t.trim() === \"2.\") {\n    idxTwo = i;\n  }\n}\nif (idxJava === -1) {\n  throw new Error(\"Could not re-locate 'In .java files attached' paragraph\");\n}\n\nconst javaOoxml = flatOpc(\n  \"<w:p>\" +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>In .java</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> files attached</w:t></w:r>' +\n    \"</w:p>\" +\n    \"<w:p>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"</w:p>\"\n);\n\nconst javaParagraph = paragraphs.items[idxJava];\n// Expand the replace-range to cover everything from the start of the\n// \"In .java files attached\" paragraph through the end of the last\n// trailing numbered paragraph (if present), so the \"1.\"/\"2.\" paragraphs\n// get removed as part of the same replace.\nlet endParagraph = javaParagraph;\nif (idxTwo !== -1) {\n  endParagraph = paragraphs.items[idxTwo];\n} else if (idxOne !== -1) {\n  endParagraph = paragraphs.items[idxOne];\n}\nconst rangeToReplace = javaParagraph\n  .getRange(\"Start\")\n  .expandTo(endParagraph.getRange(\"End\"));\nrangeToReplace.insertOoxml(javaOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same change as the target diff:\n#  1. Splits the \"Because the first catch...\" run so that\n#     \"ArithmeticException\" is wrapped in proofErr spellStart/spellEnd.\n#  2. Splits the \"B,D,A,C\" run so that \",D,A,C\" is wrapped in\n#     proofErr gramStart/gramEnd (keeping \"B\" as its own run).\n#  3. Splits \"In .java files attached\" into \"In .java\" (wrapped in\n#     proofErr gramStart/gramEnd) + \" files attached\", moves the\n#     \"_GoBack\" bookmark into its own (new, empty) paragraph right\n#     after it, and removes the trailing \"       1.\" / \"       2.\"\n#     paragraphs entirely.\n\n$d = $word.ActiveDocument\n\n$pkgOpen = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$listPPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>'\n\n# ---------------------------------------------------------------------\n# Locate the paragraphs we need to touch by their current text (robust\n# to their exact index).\n# ---------------------------------------------------------------------\n$idxException = -1\n$idxBdac = -1\n$idxJava = -1\n$idxOne = -1\n$idxTwo = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.TrimEnd()\n    if ($t.Contains(\"second catch statement with ArithmeticException\")) {\n        $idxException = $i\n    } elseif ($t -eq \"B,D,A,C\") {\n        $idxBdac = $i\n    } elseif ($t -eq \"In .java files attached\") {\n        $idxJava = $i\n    } elseif ($t -eq \"1.\") {\n        $idxOne = $i\n    } elseif ($t -eq \"2.\") {\n        $idxTwo = $i\n    }\n}\n\nif ($idxException -eq -1 -or $idxBdac -eq -1 -or $idxJava -eq -1) {\n    throw \"Could not locate expected paragraphs (exception=$idxException, bdac=$idxBdac, java=$idxJava)\"\n}\n\n# ---------------------------------------------------------------------\n# 1) \"Because the first catch statement ... ArithmeticException will\n#    never be reached.\" -> split into 3 runs with a spellStart/spellEnd\n#    proofErr pair wrapping \"ArithmeticException\".\n# ---------------------------------------------------------------------\n$exceptionXml = $pkgOpen + '<w:p>' + $listPPr + `\n    '<w:r><w:t xml:space=\"preserve\">Because the first catch statement catches the type Exception, all exceptions will be caught. Because of this, the second catch statement with </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:t>ArithmeticException</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '<w:r><w:t xml:space=\"preserve\"> will never be reached.</w:t></w:r>' + `\n    '</w:p>' + $pkgClose\n\n$d.Paragraphs($idxException).Range.InsertXML($exceptionXml)\n\n# ---------------------------------------------------------------------\n# 2) \"B,D,A,C\" -> \"B\" + proofErr gramStart/gramEnd wrapping \",D,A,C\".\n# ---------------------------------------------------------------------\n$bdacXml = $pkgOpen + '<w:p>' + $listPPr + `\n    '<w:r><w:t>B</w:t></w:r>' + `\n    '<w:proofErr w:type=\"gramStart\"/>' + `\n    '<w:r><w:t>,D,A,C</w:t></w:r>' + `\n    '<w:proofErr w:type=\"gramEnd\"/>' + `\n    '</w:p>' + $pkgClose\n\n$d.Paragraphs($idxBdac).Range.InsertXML($bdacXml)\n\n# ---------------------------------------------------------------------\n# 3) \"In .java files attached\" (+ trailing _GoBack bookmark, + the\n#    \"       1.\" / \"       2.\" paragraphs that follow) all become:\n#      <p> proofErr gramStart, \"In .java\", proofErr gramEnd, \" files\n#          attached\" </p>\n#      <p> bookmarkStart/_GoBack, bookmarkEnd </p>\n#    i.e. the bookmark moves onto its own empty paragraph and the\n#    \"1.\"/\"2.\" paragraphs disappear.\n# ---------------------------------------------------------------------\n$javaXml = $pkgOpen + `\n    '<w:p>' + `\n        '<w:proofErr w:type=\"gramStart\"/>' + `\n        '<w:r><w:t>In .java</w:t></w:r>' + `\n        '<w:proofErr w:type=\"gramEnd\"/>' + `\n        '<w:r><w:t xml:space=\"preserve\"> files attached</w:t></w:r>' + `\n    '</w:p>' + `\n    '<w:p>' + `\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' + `\n        '<w:bookmarkEnd w:id=\"0\"/>' + `\n    '</w:p>' + `\n    $pkgClose\n\n$startRange = $d.Paragraphs($idxJava).Range\n$endIdx = $idxJava\nif ($idxTwo -ne -1) {\n    $endIdx = $idxTwo\n} elseif ($idxOne -ne -1) {\n    $endIdx = $idxOne\n}\n$endRange = $d.Paragraphs($endIdx).Range\n\n$replaceRange = $d.Range($startRange.Start, $endRange.End)\n$replaceRange.InsertXML($javaXml)\n"}
